$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.465.04'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.819.96'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.75'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5091'
$ws.Range('E7').Value = '  -4.77%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3951'
$ws.Range('E8').Value = '  -1.39%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08275'
$ws.Range('E9').Value = '  +8.30%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.111'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.55'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.314'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.02'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.002'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').Value = '1.814.92'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001147'
$ws.Range('E17').Value = '  +6.63%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '92.56'
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06650'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.76'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.121'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').Value = '28.475.09'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.45'
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.272'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '21.30'
$ws.Range('E26').Value = '  +2.86%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '155.94'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('D28').Value = '2.024.05'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.411'
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.66'
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1095'
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.770'
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07054'
$ws.Range('E35').Value = '  -4.83%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2226'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02339'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.252'
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.874'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6300'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.402'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.56'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.734'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5926'
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '125.22'
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.983'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06891'
$ws.Range('E51').Value = '  +0.03%  '
